# RMS.xlsx update: add Sprint entry #9 (LoginAdapter / ValidateUser work)
# Commit message: "wrote a code in RMS Dataaccess layer - Made a class
# LoginAdapter and made a method in that class named ValidateUser and
# wrote a code for checking the login username and password"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 13: new user story entry ---
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "26-09-14"
$ws.Range("C13").Value = "study and done implementati-"
$ws.Range("D13").Value = "cmt"
$ws.Range("E13").Value = "shilpa"
$ws.Range("F13").Value = 3

# --- Row 14 & 15: continuation of the wrapped user story text ---
$ws.Range("C14").Value = "on of stored procedure and "
$ws.Range("C15").Value = "made a class in business layer"

# --- Comments for row 13 & 14 (entered after the user story column) ---
$ws.Range("G13").Value = "coding on the reset and login click is still left "

# Resize column G to fit the newly entered comment text (best-fit column width)
$ws.Columns.Item(7).AutoFit() | Out-Null

$ws.Range("G14").Value = "will make a class in data acess layer and connect it with the  web page"

# Scroll the view down and select the last comment cell, as left by the author
$ws.Range("G15").Select()
